# Cross-browser parallel suite: add Sheet2 with per-browser test data,
# fix DDT bugs on Sheet1 (column widths / selection).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 fixes -----------------------------------------------------
# Selection becomes the full data block instead of the old A3 cell.
$ws1.Range("A1:C3").Select()

# Column A no longer needs an explicit width; B/C get new (narrower/taller)
# auto-fit style widths.
$ws1.Columns.Item(2).ColumnWidth = 13.05
$ws1.Columns.Item(3).ColumnWidth = 15.4

# --- Add Sheet2 (cross browser / DDT data) -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "id"
$ws2.Range("B1").Value = "url"
$ws2.Range("C1").Value = "Expected_title"

$data = @(
  @("i1", "https://www.google.com/", "Google"),
  @("i2", "https://www.wikipedia.org/", "Wikipedia"),
  @("i3", "https://www.yahoo.com ", "ee"),
  @("i4", "https://www.gmail.com", "Gmail"),
  @("i5", "https://www.orkut.com", "orkut"),
  @("i6", "https://aot.edu.in", "AOT"),
  @("i7", "https://www.facebook.com", "Facebook")
)

$row = 2
foreach ($item in $data) {
  $ws2.Cells.Item($row, 1).Value = $item[0]
  $ws2.Cells.Item($row, 2).Value = $item[1]
  $ws2.Cells.Item($row, 3).Value = $item[2]
  $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 2), $item[1]) | Out-Null
  $row++
}
# Re-apply the shared "Hyperlink" cell style (Hyperlinks.Add creates its own
# duplicate style, so reset explicitly to reuse the workbook's one style).
$ws2.Range("B2:B8").Style = "Hyperlink"

$ws2.Columns.Item(2).ColumnWidth = 57.3

# Final selection / active cell on the new sheet.
$ws2.Range("C8").Select()

Write-Host "edit applied"
